# Swap the data between row 20 and row 21 for the columns that differ:
# A, B, D, E, F, G, H, Q, R
# (Columns C, I, P, S, ... remain identical between the two rows, so no
# change is needed there.)
#
# Note: reading via the .Value property getter on this runtime returns a
# description string instead of the actual cell value, so .Value2 is used
# for reading while .Value is used for writing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $cell20 = $ws.Range($col + "20")
    $cell21 = $ws.Range($col + "21")

    $v20 = $cell20.Value2
    $v21 = $cell21.Value2

    $cell20.Value = $v21
    $cell21.Value = $v20
}
